$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").ClearContents()
$ws.Columns.Item(1).Delete()
$ws.Range("L2").Value = "Test passed Successfully"

# Copy formatting of row 2 (A2:J2) down into row 3
$ws.Range("A2:J2").Copy()
$ws.Range("A3:J3").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A3").Value = "rohit"
$ws.Range("B3").Value = "rana"
$ws.Range("C3").Value = "Male"
$ws.Range("D3").Value = "6"
$ws.Range("E3").Value = 44278
$ws.Range("F3").Value = "Manual Tester"
$ws.Range("G3").Value = "QTP"
$ws.Range("H3").Value = "Africa"
$ws.Range("I3").Value = "Navigation Commands"
$ws.Range("J3").Value = "/src/main/resources/Feb bill.pdf"

$ws.Range("L2").Select()
